$d = $word.ActiveDocument

# The sentence about the banner-click revenue formula contained a
# misplaced "c" (cent sign): "0,3+2/50c bis 1+15/50" should read
# "0,3c+2c/50 bis 1c+15c/50". Reproduce this as a sequence of small,
# targeted in-place edits (tracked, then accepted) so the resulting
# run layout matches how Word splits runs around each edited span
# instead of collapsing the whole sentence into one rewritten run.

$d.TrackRevisions = $true

$rng = $d.Content
$found = $rng.Find.Execute(
    "0,3+2/50c bis 1+15/50",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

if ($found) {
    $start = $rng.Start

    # Apply edits from right to left so earlier offsets stay valid.
    $d.Range($start + 18, $start + 18).InsertBefore("c")   # ...+15[c]/50
    $d.Range($start + 15, $start + 15).InsertBefore("c")   # ...bis 1[c]+15
    $d.Range($start + 5, $start + 9).Text = "c/50"          # /50c -> c/50
    $d.Range($start + 3, $start + 3).InsertBefore("c")      # 0,3[c]+2
}

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
